$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell K1, mirroring the style of the other header cells (row 1)
$ws.Range("K1").Value = "PRODUCTO"
$ws.Cells.Item(1, 11).Font.Bold = $true
$ws.Cells.Item(1, 11).HorizontalAlignment = -4108

# Fill K2:K216 with "SORGO"
$lastRow = 216
$ws.Range("K2:K" + $lastRow).Value = "SORGO"
